$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before K — shifts old columns K:AB (addSequence..sequence)
# right by one, landing on L:AC, matching the diff's per-column shift.
$ws.Columns("K:K").Insert() | Out-Null

# The old "Click on delete button in confirm delete message" test description is
# replaced by the new test description for this testcase (genetic feature /
# magnifying-glass alignment check).
$ws.Range("B2").Value = "Check vertically align magnifying glass on literature evidence tab details and trait components for genetic feature"

# New column K holds the addGFSymbol / AAP55168 pair introduced by this change.
$ws.Range("K1").Value = "addGFSymbol"
$ws.Range("K2").Value = "AAP55168"

# Column widths for the newly-relevant columns (F, I, J:K, L) plus the width that
# moved from old column L (12) to new column M (13, unchanged at 44).
$ws.Columns("F:F").ColumnWidth = 18.666666666666664
$ws.Columns("I:I").ColumnWidth = 15.833333333333332
$ws.Columns("J:K").ColumnWidth = 24
$ws.Columns("L:L").ColumnWidth = 19.5

# Row 2 no longer maxes out at 409.5pt — it's now an explicit custom height.
$ws.Rows("2:2").RowHeight = 288

# Update the view's active selection to the new location.
$ws.Range("H2").Select() | Out-Null

Write-Output "edit applied"
